$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.201.97"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "3.147.58"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.82"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.37"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.139.87"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.19"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "3.668.07"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.23"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "64.028.18"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").Value = "3.146.23"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.21"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.48"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.04"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("E25").Value = "  +6.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.19"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  +8.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.33"
$ws.Range("E31").Value = "  +6.59%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").Value = "0.0₃0829"
$ws.Range("E35").Value = "  -5.71%  "
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.18"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "461.72"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.31"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.23"
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.294"
$ws.Range("E43").Value = "  +4.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0373"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "2.921.64"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.60"
$ws.Range("E46").Value = "  +8.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.108"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.93"
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("E51").Value = "  -0.81%  "
